$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "languages"
$ws.Range("K6").Select()
